# Update "想去人数" (F column) counts with refreshed data.
# Mirrors the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 3061
$ws.Range("F9").Value  = 7199
$ws.Range("F17").Value = 1894
$ws.Range("F18").Value = 1754
$ws.Range("F24").Value = 1321
$ws.Range("F31").Value = 51
$ws.Range("F32").Value = 2599
$ws.Range("F33").Value = 2897
$ws.Range("F34").Value = 2127
$ws.Range("F35").Value = 90
$ws.Range("F42").Value = 358
$ws.Range("F46").Value = 109

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value  = 78
$ws.Range("F17").Value = 510

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value  = 556
$ws.Range("F6").Value  = 1791
$ws.Range("F8").Value  = 2850
$ws.Range("F10").Value = 1038
$ws.Range("F14").Value = 7914

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 3061
$ws.Range("F6").Value  = 1791
$ws.Range("F8").Value  = 2850
$ws.Range("F9").Value  = 7199
$ws.Range("F18").Value = 1894
$ws.Range("F19").Value = 78
$ws.Range("F23").Value = 1321
$ws.Range("F31").Value = 510
$ws.Range("F34").Value = 51
$ws.Range("F35").Value = 2599
$ws.Range("F36").Value = 2897
$ws.Range("F37").Value = 2127
$ws.Range("F38").Value = 90
$ws.Range("F43").Value = 358
